$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 6).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $val = $cell.Value()
    if ($val -ne $null) {
        $trimmed = $val -replace '^\s+', ''
        if ($trimmed -ne $val) {
            $cell.Value = $trimmed
        }
    }
}
